# Exercise 4: rename Sheet1 -> OrderDetailsData, add YearlyIncome and
# MonthlyIncome summary tabs with bordered/currency-formatted tables.

$wb = $excel.ActiveWorkbook
$orderDetails = $wb.Worksheets.Item(1)
$orderDetails.Name = "OrderDetailsData"

# Tidy the active sheet's view (drop the frozen/scrolled selection state).
$orderDetails.Activate()
$orderDetails.Range("A1").Select()

$currencyFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# ---------------------------------------------------------------------
# YearlyIncome
# ---------------------------------------------------------------------
$yearly = $wb.Worksheets.Add([Type]::Missing, $orderDetails)
$yearly.Name = "YearlyIncome"

$yearly.Range("A1").Value = "Year"
$yearly.Range("B1").Value = 2016
$yearly.Range("C1").Value = 2017
$yearly.Range("D1").Value = 2018
$yearly.Range("A1:D1").Font.Bold = $true

$yearly.Range("A2").Value = "Total "

$yearly.Range("A1:D2").Borders.LineStyle = 1
$yearly.Range("B2:D2").NumberFormat = $currencyFormat

$yearly.Columns.Item(1).ColumnWidth = 10
$yearly.Columns.Item(2).ColumnWidth = 13.570312
$yearly.Columns.Item(3).ColumnWidth = 17.425781
$yearly.Columns.Item(4).ColumnWidth = 17

# ---------------------------------------------------------------------
# MonthlyIncome
# ---------------------------------------------------------------------
$monthly = $wb.Worksheets.Add([Type]::Missing, $yearly)
$monthly.Name = "MonthlyIncome"

$monthly.Range("A1").Value = "Month"
$monthly.Range("B1").Value = 2016
$monthly.Range("C1").Value = 2017
$monthly.Range("D1").Value = 2018
$monthly.Range("A1:D1").Font.Bold = $true

for ($m = 1; $m -le 12; $m++) {
    $monthly.Cells.Item($m + 1, 1).Value = $m
}

$monthly.Range("A14").Value = "Total"
$monthly.Range("A14:D14").Font.Bold = $true

$monthly.Range("A1:D14").Borders.LineStyle = 1
$monthly.Range("B2:D14").NumberFormat = $currencyFormat

$monthly.Columns.Item(1).ColumnWidth = 10
$monthly.Columns.Item(2).ColumnWidth = 13.570312
$monthly.Columns.Item(3).ColumnWidth = 17.425781
$monthly.Columns.Item(4).ColumnWidth = 17

$orderDetails.Activate()
